$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# After the PCB check, the "L1" BOM line (2.2 uH Coilcraft inductor, row 7)
# was removed from the BOM, and everything below it shifted up one row.

# Helper: (re)create a hyperlink on a cell without Excel clobbering the
# cell's existing number format / border / fill with the default
# "Hyperlink" style - we stash the original formatting on a scratch cell,
# add the hyperlink, then paste the formatting back.
function Add-HyperlinkKeepFormat($ws, $addr, $target, $showDisplay) {
    $scratch = $ws.Range("ZZ1")
    $ws.Range($addr).Copy($scratch)
    if ($showDisplay) {
        $ws.Hyperlinks.Add($ws.Range($addr), $target, "", "", $target) | Out-Null
    } else {
        $ws.Hyperlinks.Add($ws.Range($addr), $target) | Out-Null
    }
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $scratch.Clear()
}

# Cells whose hyperlink "display" text is explicitly cached (differs from
# just relying on the cell's own text), based on the original workbook.
$cellsWithExplicitDisplay = @("K8", "D9", "K6")

# 1. Capture all existing hyperlinks (address, target URL, whether they
#    need an explicit display string) before anything moves.
$hlinks = @()
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address($false, $false)
    $hlinks += [PSCustomObject]@{
        Address     = $addr
        Target      = $hl.Address
        ShowDisplay = $cellsWithExplicitDisplay -contains $addr
    }
}

# 2. Remove all hyperlinks so stale ones don't linger on the wrong cells
#    once the row shifts.
$ws.Hyperlinks.Delete()

# 3. Delete the whole row - cells below shift up to fill the gap.
$ws.Rows(7).Delete()

# 4. Recreate the surviving hyperlinks (skip the one that lived on the
#    deleted L1 row, K7, which pointed at the Coilcraft part page) at
#    their shifted addresses.
foreach ($h in $hlinks) {
    if ($h.Address -eq "K7") {
        continue
    }

    $col = $h.Address -replace '[0-9]+', ''
    $row = [int]($h.Address -replace '[A-Z]+', '')
    if ($row -gt 7) {
        $row = $row - 1
    }
    $newAddr = "$col$row"

    Add-HyperlinkKeepFormat $ws $newAddr $h.Target $h.ShowDisplay
}

# 5. Leave the selection where the author ended up after the edit.
$ws.Range("F21").Select() | Out-Null
